# "Generate Report for Handback"
#
# The localization-status report tracks the handoff/handback lifecycle of
# each source file, per target locale. The file
# "fabd5837-7993-46a6-a1bb-85f8990fe1cc.md" has just been handed back (in
# sync with en-US) for both the zh-cn and de-de locales, so:
#   * its Status flips from "Ready for handoff" to
#     "Handed back: in sync with en-US" everywhere it is reported
#     (Overview summary sheet + each per-locale detail sheet), and
#   * the per-locale "Latest Handback DateTime" timestamp is stamped with
#     the moment the handback report was generated.

$wb = $excel.ActiveWorkbook

$newStatus = "Handed back: in sync with en-US"

# --- Overview sheet: one status column per locale -------------------------
$overview = $wb.Sheets.Item("Overview")
$overview.Range("B3").Value = $newStatus   # zh-cn status
$overview.Range("C3").Value = $newStatus   # de-de status

# --- zh-cn detail sheet -----------------------------------------------------
$zhcn = $wb.Sheets.Item("zh-cn")
$zhcn.Range("B3").Value = $newStatus                    # Status
$zhcn.Range("G3").Value = "2016-03-03 15:11:28"          # Latest Handback DateTime

# --- de-de detail sheet -----------------------------------------------------
$dede = $wb.Sheets.Item("de-de")
$dede.Range("B3").Value = $newStatus                    # Status
$dede.Range("G3").Value = "2016-03-03 15:11:54"          # Latest Handback DateTime
